# Apply the three content changes from the diff to the "Metadata" sheet
# (sheet1.xml / xl/worksheets/sheet1.xml):
#   - B7  (row "Experimental"): was empty -> "false"
#   - B8  (row "Date"): "2025-11-28T14:35:57+00:00" -> "2025-11-30T13:08:37+00:00"
#   - B17 (row "Description"): was empty -> "Methods for VO2max estimation and measurement"
# Sheet2 ("Concepts") content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Leading apostrophe forces Excel to store this as literal text ("false")
# rather than auto-converting the word "false" into a Boolean value.
$ws.Range("B7").Value = "'false"

$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

$ws.Range("B17").Value = "Methods for VO2max estimation and measurement"
